# Apply the 2025-11-17 12:37 JST scrape refresh to the "ランサーズ" sheet (ActiveSheet).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row-by-row cell values (header row 1 is untouched) ---
# Row 2: 【Next.js × TypeScript × Tailwi
$ws.Cells.Item(2, 1).Value = "2025-11-17 12:37:29"
$ws.Cells.Item(2, 2).Value = "【Next.js × TypeScript × Tailwind】コンポーネント制作パートナー募集!"
$ws.Cells.Item(2, 3).Value = "システム開発"
$ws.Cells.Item(2, 4).Value = "100,000 円 ~ 200,000 円 / 固定"
$ws.Cells.Item(2, 5).Value = "期限情報なし"
$ws.Cells.Item(2, 6).Value = "https://www.lancers.jp/work/detail/5428507"
$ws.Cells.Item(2, 7).Value = 528
$ws.Cells.Item(2, 8).Value = "🔥AI,Next.js"

# Row 3: 大企業の業務効率化AIプロジェクトの技術方針策定を支援するA
$ws.Cells.Item(3, 1).Value = "2025-11-17 12:37:29"
$ws.Cells.Item(3, 2).Value = "大企業の業務効率化AIプロジェクトの技術方針策定を支援するAIテックリード募集"
$ws.Cells.Item(3, 3).Value = "システム開発"
$ws.Cells.Item(3, 4).Value = "300,000 円 ~ 500,000 円 / 固定"
$ws.Cells.Item(3, 5).Value = "期限情報なし"
$ws.Cells.Item(3, 6).Value = "https://www.lancers.jp/work/detail/5423720"
$ws.Cells.Item(3, 7).Value = 385
$ws.Cells.Item(3, 8).Value = "🔥AI,Ai ◆効率化"

# Row 4: Stable Diffusionに詳しいLoRAなどを用いた
$ws.Cells.Item(4, 1).Value = "2025-11-17 12:37:29"
$ws.Cells.Item(4, 2).Value = "Stable Diffusionに詳しいLoRAなどを用いた画像生成AIエンジニア募集"
$ws.Cells.Item(4, 3).Value = "システム開発"
$ws.Cells.Item(4, 4).Value = "300,000 円 ~ 500,000 円 / 固定"
$ws.Cells.Item(4, 5).Value = "期限情報なし"
$ws.Cells.Item(4, 6).Value = "https://www.lancers.jp/work/detail/5416328"
$ws.Cells.Item(4, 7).Value = 310
$ws.Cells.Item(4, 8).Value = "🔥AI,Ai"

# Row 5: 製造業向けAI戦略アドバイザー募集(事業価値試算・プロジェク
$ws.Cells.Item(5, 1).Value = "2025-11-17 12:37:29"
$ws.Cells.Item(5, 2).Value = "製造業向けAI戦略アドバイザー募集(事業価値試算・プロジェクト推進支援)"
$ws.Cells.Item(5, 3).Value = "システム開発"
$ws.Cells.Item(5, 4).Value = "20,000 円 ~ 50,000 円 / 固定"
$ws.Cells.Item(5, 5).Value = "期限情報なし"
$ws.Cells.Item(5, 6).Value = "https://www.lancers.jp/work/detail/5419380"
$ws.Cells.Item(5, 7).Value = 298
$ws.Cells.Item(5, 8).Value = "🔥AI,Ai"

# Row 6: 医療系機械学習モデル活用のGUIアプリ開発
$ws.Cells.Item(6, 1).Value = "2025-11-17 12:37:29"
$ws.Cells.Item(6, 2).Value = "医療系機械学習モデル活用のGUIアプリ開発"
$ws.Cells.Item(6, 3).Value = "システム開発"
$ws.Cells.Item(6, 4).Value = "300,000 円 ~ 500,000 円 / 固定"
$ws.Cells.Item(6, 5).Value = "期限情報なし"
$ws.Cells.Item(6, 6).Value = "https://www.lancers.jp/work/detail/5435875"
$ws.Cells.Item(6, 7).Value = 225
$ws.Cells.Item(6, 8).Value = "🔥機械学習 ◆開発 ◇アプリ"

# Row 7: 【自動運転プロジェクト経験者募集】実証実験・開発を推進するプ
$ws.Cells.Item(7, 1).Value = "2025-11-17 12:37:29"
$ws.Cells.Item(7, 2).Value = "【自動運転プロジェクト経験者募集】実証実験・開発を推進するプロジェクトマネージャー"
$ws.Cells.Item(7, 3).Value = "システム開発"
$ws.Cells.Item(7, 4).Value = "200,000 円 ~ 300,000 円 / 固定"
$ws.Cells.Item(7, 5).Value = "期限情報なし"
$ws.Cells.Item(7, 6).Value = "https://www.lancers.jp/work/detail/5431107"
$ws.Cells.Item(7, 7).Value = 68
$ws.Cells.Item(7, 8).Value = "◆開発"

# Row 8: UTAGE構築代行|ヒアリングから構築までお任せしたいです。
$ws.Cells.Item(8, 1).Value = "2025-11-17 12:37:29"
$ws.Cells.Item(8, 2).Value = "UTAGE構築代行|ヒアリングから構築までお任せしたいです。"
$ws.Cells.Item(8, 3).Value = "システム開発"
$ws.Cells.Item(8, 4).Value = "50,000 円 ~ 100,000 円 / 固定"
$ws.Cells.Item(8, 5).Value = "期限情報なし"
$ws.Cells.Item(8, 6).Value = "https://www.lancers.jp/work/detail/5429882"
$ws.Cells.Item(8, 7).Value = 18
$ws.Cells.Item(8, 8).ClearContents()

# Row 9: 【急募】Wartalesの武器アイコンとモデルを日本刀に差し
$ws.Cells.Item(9, 1).Value = "2025-11-17 12:37:29"
$ws.Cells.Item(9, 2).Value = "【急募】Wartalesの武器アイコンとモデルを日本刀に差し替え"
$ws.Cells.Item(9, 3).Value = "システム開発"
$ws.Cells.Item(9, 4).Value = "20,000 円 ~ 30,000 円 / 募集期間 5 日、取引期間 0 日"
$ws.Cells.Item(9, 5).Value = "期限情報なし"
$ws.Cells.Item(9, 6).Value = "https://www.lancers.jp/work/detail/5435667"
$ws.Cells.Item(9, 7).Value = 10
$ws.Cells.Item(9, 8).ClearContents()

# Row 10: ロリポップ!レンタルサーバーの不具合を解決したい
$ws.Cells.Item(10, 1).Value = "2025-11-17 12:37:29"
$ws.Cells.Item(10, 2).Value = "ロリポップ!レンタルサーバーの不具合を解決したい"
$ws.Cells.Item(10, 3).Value = "システム開発"
$ws.Cells.Item(10, 4).Value = "5,000 円 ~ 10,000 円 / 固定"
$ws.Cells.Item(10, 5).Value = "期限情報なし"
$ws.Cells.Item(10, 6).Value = "https://www.lancers.jp/work/detail/5435519"
$ws.Cells.Item(10, 7).Value = 10
$ws.Cells.Item(10, 8).ClearContents()

# --- Rebuild hyperlinks on the URL column (F) in row order so relationship ids line up 1:1 ---
$ws.Cells.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("F2"), "https://www.lancers.jp/work/detail/5428507")
$ws.Hyperlinks.Add($ws.Range("F3"), "https://www.lancers.jp/work/detail/5423720")
$ws.Hyperlinks.Add($ws.Range("F4"), "https://www.lancers.jp/work/detail/5416328")
$ws.Hyperlinks.Add($ws.Range("F5"), "https://www.lancers.jp/work/detail/5419380")
$ws.Hyperlinks.Add($ws.Range("F6"), "https://www.lancers.jp/work/detail/5435875")
$ws.Hyperlinks.Add($ws.Range("F7"), "https://www.lancers.jp/work/detail/5431107")
$ws.Hyperlinks.Add($ws.Range("F8"), "https://www.lancers.jp/work/detail/5429882")
$ws.Hyperlinks.Add($ws.Range("F9"), "https://www.lancers.jp/work/detail/5435667")
$ws.Hyperlinks.Add($ws.Range("F10"), "https://www.lancers.jp/work/detail/5435519")
$ws.Range("F2:F10").Style = "Hyperlink"

# --- Column width updates (D widened for longer price strings, H widened for skill tags) ---
$ws.Columns.Item(4).ColumnWidth = 40.166666666666664
$ws.Columns.Item(8).ColumnWidth = 15.166666666666666
